$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.717.98'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '3.108.83'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '231.52'
$ws.Range('E5').Value = '  +5.16%  '
$ws.Range('D6').Value = '627.29'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '1.09'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').Value = '0.366'
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '2.880.65'
$ws.Range('E10').Value = '  -9.53%  '
$ws.Range('D11').Value = '0.722'
$ws.Range('E11').Value = '  -4.18%  '
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D13').Value = '36.56'
$ws.Range('E13').Value = '  +3.52%  '
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = '5.49'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').Value = '90.683.84'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '3.693.40'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').Value = '3.128.41'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').Value = '3.78'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').Value = '14.10'
$ws.Range('E20').Value = '  -2.14%  '
$ws.Range('D21').Value = '0.0000208'
$ws.Range('E21').Value = '  -5.91%  '
$ws.Range('D22').Value = '440.09'
$ws.Range('E22').Value = '  -1.02%  '
$ws.Range('D23').Value = '5.54'
$ws.Range('E23').Value = '  +6.26%  '
$ws.Range('D24').Value = '8.89'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').Value = '5.68'
$ws.Range('E25').Value = '  -6.63%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '89.00'
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '12.37'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = '3.286.57'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '9.38'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').Value = '0.160'
$ws.Range('E31').Value = '  -3.05%  '
$ws.Range('D32').Value = '0.197'
$ws.Range('E32').Value = '  +17.13%  '
$ws.Range('D33').Value = '26.20'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('E34').Value = '  -8.40%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.151'
$ws.Range('E35').Value = '  +2.81%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '509.19'
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').Value = '3.76'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').Value = '1.93'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '7.03'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').Value = '1.28'
$ws.Range('E40').Value = '  -2.97%  '
$ws.Range('D41').Value = '0.411'
$ws.Range('E41').Value = '  -1.29%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0864'
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').Value = '22.18'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '3.34'
$ws.Range('E45').Value = '  +48.80%  '
$ws.Range('D46').Value = '1.90'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').Value = '150.39'
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '45.38'
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '0.685'
$ws.Range('E49').Value = '  +4.37%  '
$ws.Range('D50').Value = '1.34'
$ws.Range('E50').Value = '  -2.04%  '
$ws.Range('D51').Value = '4.43'
$ws.Range('E51').Value = '  +0.34%  '
